$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Compare-Ordinal($s1, $s2) {
    $len1 = $s1.Length
    $len2 = $s2.Length
    $minLen = $len1
    if ($len2 -lt $minLen) { $minLen = $len2 }
    for ($k = 0; $k -lt $minLen; $k++) {
        $c1 = [int][char]$s1.Substring($k, 1)
        $c2 = [int][char]$s2.Substring($k, 1)
        if ($c1 -lt $c2) { return -1 }
        if ($c1 -gt $c2) { return 1 }
    }
    if ($len1 -lt $len2) { return -1 }
    if ($len1 -gt $len2) { return 1 }
    return 0
}

function Sort-OrdinalList($list) {
    $n = $list.Count
    for ($i = 1; $i -lt $n; $i++) {
        $key = $list[$i]
        $j = $i - 1
        while ($j -ge 0 -and (Compare-Ordinal $list[$j] $key) -gt 0) {
            $list[$j+1] = $list[$j]
            $j = $j - 1
        }
        $list[$j+1] = $key
    }
    return $list
}

$lastRow = $ws.UsedRange.Rows.Count
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text
    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ", "
        if ($parts.Count -gt 1) {
            $sortedParts = Sort-OrdinalList $parts
            $newText = $sortedParts -join ", "
            if ($newText -cne $text) {
                $cell.Value = $newText
            }
        }
    }
}
